$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Atomic Habits): mark as back "In Stock" and clear "Who Checked" since the book was returned.
$ws.Range("F4").Value = "Yes"
$ws.Range("G4").Value = ""
